# This workbook is a weekly price log. A new week's worth of records
# (2 rows) is inserted at the top of the data block (row 109), pushing all
# the existing data down by 2 rows. The two rows that fall off the bottom of
# the previous range become new rows 165/166, extending the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 109:110 - this shifts old rows 109..164 down to
# 111..166, carrying all their existing data and formatting with them.
$ws.Rows("109:110").Insert()

# Populate the two newly-inserted rows with the new week's data.

# Row 109
$ws.Cells.Item(109, 1).Value = 9
$ws.Cells.Item(109, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(109, 3).Value = "Metropolitana"
$ws.Cells.Item(109, 4).Value = 45089
$ws.Cells.Item(109, 5).Value = 13
$ws.Cells.Item(109, 6).Value = 100114002
$ws.Cells.Item(109, 7).Value = "Camote"
$ws.Cells.Item(109, 8).Value = "Sin especificar"
$ws.Cells.Item(109, 9).Value = "Primera"
$ws.Cells.Item(109, 10).Value = 700
$ws.Cells.Item(109, 11).Value = 18000
$ws.Cells.Item(109, 12).Value = 19000
$ws.Cells.Item(109, 13).Value = 18500
$ws.Cells.Item(109, 14).Value = "`$/caja 18 kilos"
$ws.Cells.Item(109, 15).Value = "Perú"
$ws.Cells.Item(109, 16).Value = 1028
$ws.Cells.Item(109, 17).Value = 18
$ws.Cells.Item(109, 18).Value = "Hortaliza"

# Row 110
$ws.Cells.Item(110, 1).Value = 9
$ws.Cells.Item(110, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(110, 3).Value = "Metropolitana"
$ws.Cells.Item(110, 4).Value = 45089
$ws.Cells.Item(110, 5).Value = 13
$ws.Cells.Item(110, 6).Value = 100114002
$ws.Cells.Item(110, 7).Value = "Camote"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 520
$ws.Cells.Item(110, 11).Value = 16000
$ws.Cells.Item(110, 12).Value = 17000
$ws.Cells.Item(110, 13).Value = 16500
$ws.Cells.Item(110, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(110, 15).Value = "Perú"
$ws.Cells.Item(110, 16).Value = 917
$ws.Cells.Item(110, 17).Value = 18
$ws.Cells.Item(110, 18).Value = "Hortaliza"
